# Update column G ("K") values per regenerated save_data (Strike# -> K)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 1
    6  = 2
    7  = 0
    8  = 2
    10 = 0
    11 = 2
    12 = 2
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 6
    18 = 2
    19 = 3
    20 = 0
    21 = 2
    22 = 6
    23 = 2
    24 = 0
    25 = 1
    26 = 3
    27 = 0
    28 = 1
    29 = 3
    30 = 2
    31 = 3
    32 = 1
    33 = 1
    34 = 0
    35 = 2
    36 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
